# Update hospital capacity index: convert index into numeric
# Root cause: B4 (mean_travel_time_mins) and C4 (median_travel_time_mins)
# for "Regionalverband Saarbrücken" change, which shifts the min of each
# column used for min-max scaling, so the scaled columns E and F are
# recalculated for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Root values that actually changed (row 4) ---
$ws.Range("B4").Value = 12.51
$ws.Range("C4").Value = 9.43

# --- Recalculated min-max scaled values ---

# Column E (mean_travel_time_mins_scaled)
$ws.Range("E2").Value = 0.8220155038759691
$ws.Range("E3").Value = 0.5444961240310074
$ws.Range("E5").Value = 0.7072868217054261
$ws.Range("E7").Value = 0.995968992248062

# Column F (median_travel_time_mins_scaled)
$ws.Range("F2").Value = 0.3061611374407583
$ws.Range("F3").Value = 0.7611374407582939
$ws.Range("F5").Value = 0.5924170616113744
$ws.Range("F7").Value = 0.78957345971564
